# Add "Distribution Center" (C) and "Affected" (D) columns to the
# "SehirlerBolgeler" worksheet, populated with the per-city flag data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SehirlerBolgeler")

# New header cells
$ws.Cells.Item(1, 3).Value = "Distribution Center"
$ws.Cells.Item(1, 4).Value = "Affected"

# Per-row values for column C (Distribution Center) and D (Affected),
# one entry per data row (rows 2..82, cities 1..81).
$cVals = @(0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,1,1,0,0,1,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$dVals = @(1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,1,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0)

for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
}

# Match the author's final selection (cell D2 active/selected).
$null = $ws.Range("D2").Select()
